$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: swap full row content (columns B..AC) between paired rows ---
# These pairs represent records whose data got interchanged (same fixture date,
# order of the two matches swapped), while column A (sequence id) stays fixed.
$rowPairs = @(
    @(3, 4),
    @(7, 8),
    @(65, 66),
    @(82, 84),
    @(90, 91),
    @(94, 95),
    @(98, 99),
    @(104, 105),
    @(110, 111),
    @(119, 120),
    @(130, 131)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    for ($col = 2; $col -le 29; $col++) {
        $c1 = $ws.Cells.Item($r1, $col)
        $c2 = $ws.Cells.Item($r2, $col)
        $v1 = $c1.Value2
        $v2 = $c2.Value2
        $c1.Value = $v2
        $c2.Value = $v1
    }
}

# --- Part 2: rename swapped team names in HomeTeam (F) / AwayTeam (G) cells ---
# "AE Altos" <-> "Manaus" and "Botafogo PB" <-> "Paysandu" were swapped in the
# underlying shared-string table; apply the same rename to every remaining cell
# that still references one of those four team names (rows already handled above
# are excluded).
function Swap-TeamName($name) {
    if ($name -eq "AE Altos") { return "Manaus" }
    if ($name -eq "Manaus") { return "AE Altos" }
    if ($name -eq "Botafogo PB") { return "Paysandu" }
    if ($name -eq "Paysandu") { return "Botafogo PB" }
    return $name
}

$renameCells = @(
    @(13, "G"),
    @(15, "G"),
    @(16, "G"),
    @(20, "G"),
    @(22, "G"),
    @(23, "G"),
    @(28, "G"),
    @(29, "F"),
    @(32, "F"),
    @(36, "F"),
    @(38, "F"),
    @(39, "F"),
    @(42, "F"),
    @(45, "F"),
    @(49, "F"),
    @(49, "G"),
    @(53, "G"),
    @(54, "G"),
    @(55, "F"),
    @(61, "G"),
    @(62, "F"),
    @(63, "G"),
    @(64, "G"),
    @(73, "F"),
    @(75, "F"),
    @(77, "F"),
    @(81, "F"),
    @(86, "F"),
    @(87, "F"),
    @(87, "G"),
    @(88, "F"),
    @(97, "G"),
    @(103, "G"),
    @(108, "F"),
    @(108, "G"),
    @(114, "G"),
    @(116, "F"),
    @(117, "F"),
    @(118, "G"),
    @(123, "F"),
    @(125, "F"),
    @(127, "F"),
    @(128, "F"),
    @(134, "G"),
    @(135, "G"),
    @(136, "G"),
    @(140, "G"),
    @(144, "F"),
    @(145, "F"),
    @(146, "G"),
    @(149, "G"),
    @(152, "F"),
    @(152, "G"),
    @(155, "F"),
    @(155, "G"),
    @(159, "F"),
    @(161, "F"),
    @(162, "G"),
    @(163, "G")
)

foreach ($item in $renameCells) {
    $row = $item[0]
    $colLetter = $item[1]
    $cell = $ws.Range("$colLetter$row")
    $current = $cell.Value2
    $cell.Value = Swap-TeamName $current
}

Write-Host "Done applying Brazil Serie C updates."